# Apply updated TPM-derived values (ligand/receptor/edge columns) per the new TPM data.
# Each cluster-level ligand (E-J) and receptor (K-P) block, plus the per-pair edge
# weights/specificities (Q-T) that derive from them, are rewritten with recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.532141"
$ws.Range("H2").Value = [double]"4.596423"
$ws.Range("I2").Value = [double]"0.08900664250669833"
$ws.Range("J2").Value = [double]"0.08900664250669831"
$ws.Range("M2").Value = [double]"2.733663333333333"
$ws.Range("N2").Value = [double]"8.200989999999999"
$ws.Range("O2").Value = [double]"0.04037266183309663"
$ws.Range("P2").Value = [double]"0.04037266183309663"
$ws.Range("Q2").Value = [double]"4.188357673196666"
$ws.Range("R2").Value = [double]"37.69521905876999"
$ws.Range("S2").Value = [double]"0.003593435078822256"
$ws.Range("T2").Value = [double]"0.003593435078822255"
# Row 3
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.532141"
$ws.Range("H3").Value = [double]"4.596423"
$ws.Range("I3").Value = [double]"0.08900664250669833"
$ws.Range("J3").Value = [double]"0.08900664250669831"
$ws.Range("O3").Value = [double]"0.6389522306252696"
$ws.Range("P3").Value = [double]"0.6389522306252696"
$ws.Range("Q3").Value = [double]"66.28645118840333"
$ws.Range("R3").Value = [double]"596.5780606956299"
$ws.Range("S3").Value = [double]"0.05687099277012084"
$ws.Range("T3").Value = [double]"0.05687099277012082"
# Row 4
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.532141"
$ws.Range("H4").Value = [double]"4.596423"
$ws.Range("I4").Value = [double]"0.08900664250669833"
$ws.Range("J4").Value = [double]"0.08900664250669831"
$ws.Range("M4").Value = [double]"21.46453166666667"
$ws.Range("N4").Value = [double]"64.393595"
$ws.Range("O4").Value = [double]"0.3170032929137071"
$ws.Range("P4").Value = [double]"0.317003292913707"
$ws.Range("Q4").Value = [double]"32.88668901229833"
$ws.Range("R4").Value = [double]"295.980201110685"
$ws.Range("S4").Value = [double]"0.02821539876581651"
$ws.Range("T4").Value = [double]"0.02821539876581649"
# Row 5
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"1.532141"
$ws.Range("H5").Value = [double]"4.596423"
$ws.Range("I5").Value = [double]"0.08900664250669833"
$ws.Range("J5").Value = [double]"0.08900664250669831"
$ws.Range("M5").Value = [double]"0.2486213333333333"
$ws.Range("N5").Value = [double]"0.745864"
$ws.Range("O5").Value = [double]"0.003671814627926724"
$ws.Range("P5").Value = [double]"0.003671814627926724"
$ws.Range("Q5").Value = [double]"0.3809229382746667"
$ws.Range("R5").Value = [double]"3.428306444472"
$ws.Range("S5").Value = [double]"0.0003268158919387395"
$ws.Range("T5").Value = [double]"0.0003268158919387394"
# Row 6
$ws.Range("I6").Value = [double]"0.6169137955113024"
$ws.Range("J6").Value = [double]"0.6169137955113023"
$ws.Range("M6").Value = [double]"2.733663333333333"
$ws.Range("N6").Value = [double]"8.200989999999999"
$ws.Range("O6").Value = [double]"0.04037266183309663"
$ws.Range("P6").Value = [double]"0.04037266183309663"
$ws.Range("Q6").Value = [double]"29.02991907526666"
$ws.Range("R6").Value = [double]"261.2692716774"
$ws.Range("S6").Value = [double]"0.02490645204634994"
$ws.Range("T6").Value = [double]"0.02490645204634993"
# Row 7
$ws.Range("I7").Value = [double]"0.6169137955113024"
$ws.Range("J7").Value = [double]"0.6169137955113023"
$ws.Range("O7").Value = [double]"0.6389522306252696"
$ws.Range("P7").Value = [double]"0.6389522306252696"
$ws.Range("S7").Value = [double]"0.3941784457454481"
$ws.Range("T7").Value = [double]"0.3941784457454481"
# Row 8
$ws.Range("I8").Value = [double]"0.6169137955113024"
$ws.Range("J8").Value = [double]"0.6169137955113023"
$ws.Range("M8").Value = [double]"21.46453166666667"
$ws.Range("N8").Value = [double]"64.393595"
$ws.Range("O8").Value = [double]"0.3170032929137071"
$ws.Range("P8").Value = [double]"0.317003292913707"
$ws.Range("Q8").Value = [double]"227.9408768716333"
$ws.Range("R8").Value = [double]"2051.4678918447"
$ws.Range("S8").Value = [double]"0.1955637046209762"
$ws.Range("T8").Value = [double]"0.1955637046209761"
# Row 9
$ws.Range("I9").Value = [double]"0.6169137955113024"
$ws.Range("J9").Value = [double]"0.6169137955113023"
$ws.Range("M9").Value = [double]"0.2486213333333333"
$ws.Range("N9").Value = [double]"0.745864"
$ws.Range("O9").Value = [double]"0.003671814627926724"
$ws.Range("P9").Value = [double]"0.003671814627926724"
$ws.Range("Q9").Value = [double]"2.640214359626667"
$ws.Range("R9").Value = [double]"23.76192923664"
$ws.Range("S9").Value = [double]"0.002265193098528196"
$ws.Range("T9").Value = [double]"0.002265193098528196"
# Row 10
$ws.Range("G10").Value = [double]"4.902263666666666"
$ws.Range("H10").Value = [double]"14.706791"
$ws.Range("I10").Value = [double]"0.2847871244569372"
$ws.Range("J10").Value = [double]"0.2847871244569371"
$ws.Range("M10").Value = [double]"2.733663333333333"
$ws.Range("N10").Value = [double]"8.200989999999999"
$ws.Range("O10").Value = [double]"0.04037266183309663"
$ws.Range("P10").Value = [double]"0.04037266183309663"
$ws.Range("Q10").Value = [double]"13.40113843589889"
$ws.Range("R10").Value = [double]"120.61024592309"
$ws.Range("S10").Value = [double]"0.01149761427011993"
$ws.Range("T10").Value = [double]"0.01149761427011992"
# Row 11
$ws.Range("G11").Value = [double]"4.902263666666666"
$ws.Range("H11").Value = [double]"14.706791"
$ws.Range("I11").Value = [double]"0.2847871244569372"
$ws.Range("J11").Value = [double]"0.2847871244569371"
$ws.Range("O11").Value = [double]"0.6389522306252696"
$ws.Range("P11").Value = [double]"0.6389522306252696"
$ws.Range("Q11").Value = [double]"212.0912247979678"
$ws.Range("R11").Value = [double]"1908.82102318171"
$ws.Range("S11").Value = [double]"0.1819653684251163"
$ws.Range("T11").Value = [double]"0.1819653684251162"
# Row 12
$ws.Range("G12").Value = [double]"4.902263666666666"
$ws.Range("H12").Value = [double]"14.706791"
$ws.Range("I12").Value = [double]"0.2847871244569372"
$ws.Range("J12").Value = [double]"0.2847871244569371"
$ws.Range("M12").Value = [double]"21.46453166666667"
$ws.Range("N12").Value = [double]"64.393595"
$ws.Range("O12").Value = [double]"0.3170032929137071"
$ws.Range("P12").Value = [double]"0.317003292913707"
$ws.Range("Q12").Value = [double]"105.2247937115161"
$ws.Range("R12").Value = [double]"947.023143403645"
$ws.Range("S12").Value = [double]"0.0902784562322748"
$ws.Range("T12").Value = [double]"0.09027845623227478"
# Row 13
$ws.Range("G13").Value = [double]"4.902263666666666"
$ws.Range("H13").Value = [double]"14.706791"
$ws.Range("I13").Value = [double]"0.2847871244569372"
$ws.Range("J13").Value = [double]"0.2847871244569371"
$ws.Range("M13").Value = [double]"0.2486213333333333"
$ws.Range("N13").Value = [double]"0.745864"
$ws.Range("O13").Value = [double]"0.003671814627926724"
$ws.Range("P13").Value = [double]"0.003671814627926724"
$ws.Range("Q13").Value = [double]"1.218807329158222"
$ws.Range("R13").Value = [double]"10.969265962424"
$ws.Range("S13").Value = [double]"0.00104568552942617"
$ws.Range("T13").Value = [double]"0.00104568552942617"
# Row 14
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.159958"
$ws.Range("H14").Value = [double]"0.479874"
$ws.Range("I14").Value = [double]"0.009292437525062282"
$ws.Range("J14").Value = [double]"0.009292437525062281"
$ws.Range("M14").Value = [double]"2.733663333333333"
$ws.Range("N14").Value = [double]"8.200989999999999"
$ws.Range("O14").Value = [double]"0.04037266183309663"
$ws.Range("P14").Value = [double]"0.04037266183309663"
$ws.Range("Q14").Value = [double]"0.4372713194733333"
$ws.Range("R14").Value = [double]"3.935441875259999"
$ws.Range("S14").Value = [double]"0.0003751604378045169"
$ws.Range("T14").Value = [double]"0.0003751604378045168"
# Row 15
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.159958"
$ws.Range("H15").Value = [double]"0.479874"
$ws.Range("I15").Value = [double]"0.009292437525062282"
$ws.Range("J15").Value = [double]"0.009292437525062281"
$ws.Range("O15").Value = [double]"0.6389522306252696"
$ws.Range("P15").Value = [double]"0.6389522306252696"
$ws.Range("Q15").Value = [double]"6.920412781326666"
$ws.Range("R15").Value = [double]"62.28371503194"
$ws.Range("S15").Value = [double]"0.005937423684584505"
$ws.Range("T15").Value = [double]"0.005937423684584504"
# Row 16
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.159958"
$ws.Range("H16").Value = [double]"0.479874"
$ws.Range("I16").Value = [double]"0.009292437525062282"
$ws.Range("J16").Value = [double]"0.009292437525062281"
$ws.Range("M16").Value = [double]"21.46453166666667"
$ws.Range("N16").Value = [double]"64.393595"
$ws.Range("O16").Value = [double]"0.3170032929137071"
$ws.Range("P16").Value = [double]"0.317003292913707"
$ws.Range("Q16").Value = [double]"3.433423556336667"
$ws.Range("R16").Value = [double]"30.90081200703"
$ws.Range("S16").Value = [double]"0.002945733294639642"
$ws.Range("T16").Value = [double]"0.002945733294639641"
# Row 17
$ws.Range("E17").Value = [double]"2"
$ws.Range("F17").Value = [double]"0.6666666666666666"
$ws.Range("G17").Value = [double]"0.159958"
$ws.Range("H17").Value = [double]"0.479874"
$ws.Range("I17").Value = [double]"0.009292437525062282"
$ws.Range("J17").Value = [double]"0.009292437525062281"
$ws.Range("M17").Value = [double]"0.2486213333333333"
$ws.Range("N17").Value = [double]"0.745864"
$ws.Range("O17").Value = [double]"0.003671814627926724"
$ws.Range("P17").Value = [double]"0.003671814627926724"
$ws.Range("Q17").Value = [double]"0.03976897123733333"
$ws.Range("R17").Value = [double]"0.357920741136"
$ws.Range("S17").Value = [double]"3.412010803361889E-05"
$ws.Range("T17").Value = [double]"3.412010803361889E-05"
